$d = $word.ActiveDocument

# Find the run containing "Alexander Martinez" and collapse the range to its end
$rng = $d.Content
$rng.Find.Execute("Alexander Martinez", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$insertRange = $rng.Duplicate
$insertRange.Collapse(0)  # wdCollapseEnd

$insertRange.InsertAfter(" 123321")

# Apply the same run formatting as specified in the diff
$insertRange.Style = "normaltextrun"
$insertRange.Font.Name = "Aptos"
$insertRange.Font.NameAscii = "Aptos"
$insertRange.Font.Size = 14
